$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Tff3"
$ws.Range("C2").Value = "Ackr3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.461763666666667
$ws.Range("H2").Value = 4.385291
$ws.Range("I2").Value = 0.4829359810344849
$ws.Range("J2").Value = 0.482935981034485
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 16.23967033333334
$ws.Range("N2").Value = 48.71901100000001
$ws.Range("O2").Value = 0.3412424148893533
$ws.Range("P2").Value = 0.3412424148893533
$ws.Range("Q2").Value = 23.73856005191123
$ws.Range("R2").Value = 213.6470404672011
$ws.Range("S2").Value = 0.1647982404051666
$ws.Range("T2").Value = 0.1647982404051666

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Tff3"
$ws.Range("C3").Value = "Ackr3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.461763666666667
$ws.Range("H3").Value = 4.385291
$ws.Range("I3").Value = 0.4829359810344849
$ws.Range("J3").Value = 0.482935981034485
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 23.19964933333333
$ws.Range("N3").Value = 69.59894799999999
$ws.Range("O3").Value = 0.4874916916781935
$ws.Range("P3").Value = 0.4874916916781935
$ws.Range("Q3").Value = 33.91240447487422
$ws.Range("R3").Value = 305.211640273868
$ws.Range("S3").Value = 0.235427278366769
$ws.Range("T3").Value = 0.235427278366769

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Tff3"
$ws.Range("C4").Value = "Ackr3"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.461763666666667
$ws.Range("H4").Value = 4.385291
$ws.Range("I4").Value = 0.4829359810344849
$ws.Range("J4").Value = 0.482935981034485
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 8.150515666666667
$ws.Range("N4").Value = 24.451547
$ws.Range("O4").Value = 0.1712658934324533
$ws.Range("P4").Value = 0.1712658934324533
$ws.Range("Q4").Value = 11.91412766613078
$ws.Range("R4").Value = 107.227148995177
$ws.Range("S4").Value = 0.0827104622625494
$ws.Range("T4").Value = 0.0827104622625494

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Tff3"
$ws.Range("C5").Value = "Ackr3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.4482143333333333
$ws.Range("H5").Value = 1.344643
$ws.Range("I5").Value = 0.1480805917660089
$ws.Range("J5").Value = 0.1480805917660089
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 16.23967033333334
$ws.Range("N5").Value = 48.71901100000001
$ws.Range("O5").Value = 0.3412424148893533
$ws.Range("P5").Value = 0.3412424148893533
$ws.Range("Q5").Value = 7.278853012008113
$ws.Range("R5").Value = 65.50967710807302
$ws.Range("S5").Value = 0.05053137873247736
$ws.Range("T5").Value = 0.05053137873247735

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Tff3"
$ws.Range("C6").Value = "Ackr3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.4482143333333333
$ws.Range("H6").Value = 1.344643
$ws.Range("I6").Value = 0.1480805917660089
$ws.Range("J6").Value = 0.1480805917660089
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 23.19964933333333
$ws.Range("N6").Value = 69.59894799999999
$ws.Range("O6").Value = 0.4874916916781935
$ws.Range("P6").Value = 0.4874916916781935
$ws.Range("Q6").Value = 10.39841535950711
$ws.Range("R6").Value = 93.58573823556399
$ws.Range("S6").Value = 0.07218805818471964
$ws.Range("T6").Value = 0.07218805818471963

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Tff3"
$ws.Range("C7").Value = "Ackr3"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.4482143333333333
$ws.Range("H7").Value = 1.344643
$ws.Range("I7").Value = 0.1480805917660089
$ws.Range("J7").Value = 0.1480805917660089
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 8.150515666666667
$ws.Range("N7").Value = 24.451547
$ws.Range("O7").Value = 0.1712658934324533
$ws.Range("P7").Value = 0.1712658934324533
$ws.Range("Q7").Value = 3.653177945857889
$ws.Range("R7").Value = 32.878601512721
$ws.Range("S7").Value = 0.0253611548488119
$ws.Range("T7").Value = 0.0253611548488119

# Row 8
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("B8").Value = "Tff3"
$ws.Range("C8").Value = "Ackr3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 1.116849
$ws.Range("H8").Value = 3.350547
$ws.Range("I8").Value = 0.3689834271995063
$ws.Range("J8").Value = 0.3689834271995063
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 16.23967033333334
$ws.Range("N8").Value = 48.71901100000001
$ws.Range("O8").Value = 0.3412424148893533
$ws.Range("P8").Value = 0.3412424148893533
$ws.Range("Q8").Value = 18.137259572113
$ws.Range("R8").Value = 163.235336149017
$ws.Range("S8").Value = 0.1259127957517094
$ws.Range("T8").Value = 0.1259127957517094

# Row 9
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("B9").Value = "Tff3"
$ws.Range("C9").Value = "Ackr3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 1.116849
$ws.Range("H9").Value = 3.350547
$ws.Range("I9").Value = 0.3689834271995063
$ws.Range("J9").Value = 0.3689834271995063
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 23.19964933333333
$ws.Range("N9").Value = 69.59894799999999
$ws.Range("O9").Value = 0.4874916916781935
$ws.Range("P9").Value = 0.4874916916781935
$ws.Range("Q9").Value = 25.910505158284
$ws.Range("R9").Value = 233.194546424556
$ws.Range("S9").Value = 0.1798763551267049
$ws.Range("T9").Value = 0.1798763551267049

# Row 10
$ws.Range("A10").Value = "Resolving-Mac"
$ws.Range("B10").Value = "Tff3"
$ws.Range("C10").Value = "Ackr3"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 1.116849
$ws.Range("H10").Value = 3.350547
$ws.Range("I10").Value = 0.3689834271995063
$ws.Range("J10").Value = 0.3689834271995063
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 8.150515666666667
$ws.Range("N10").Value = 24.451547
$ws.Range("O10").Value = 0.1712658934324533
$ws.Range("P10").Value = 0.1712658934324533
$ws.Range("Q10").Value = 9.102895271801
$ws.Range("R10").Value = 81.926057446209
$ws.Range("S10").Value = 0.06319427632109204
$ws.Range("T10").Value = 0.06319427632109204
